$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-8 from 45243 to 45244
foreach ($r in 2..8) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45243) {
        $cell.Value = 45244
    }
}
